$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 341.8
$ws.Range("I80").Value = 382
$ws.Range("J80").Value = 324.57144
$ws.Range("K80").Value = 1146
$ws.Range("L80").Value = 973.71432
$ws.Range("M80").Value = -148
$ws.Range("N80").Value = -2969.71432
$ws.Range("H83").Value = 341.8
$ws.Range("I83").Value = 382
$ws.Range("J83").Value = 324.57144
$ws.Range("K83").Value = 3438
$ws.Range("L83").Value = 2921.14296
$ws.Range("M83").Value = 1554
$ws.Range("N83").Value = -12905.14296
$ws.Range("H88").Value = 2354.4
$ws.Range("I88").Value = 1350
$ws.Range("J88").Value = 2605.5
$ws.Range("K88").Value = 1350
$ws.Range("L88").Value = 2605.5
$ws.Range("M88").Value = -944
$ws.Range("N88").Value = -3417.5
$ws.Range("H91").Value = 2354.4
$ws.Range("I91").Value = 1350
$ws.Range("J91").Value = 2605.5
$ws.Range("K91").Value = 1350
$ws.Range("L91").Value = 2605.5
$ws.Range("M91").Value = 54
$ws.Range("N91").Value = -5413.5
$ws.Range("H92").Value = 285.44
$ws.Range("I92").Value = 260.05
$ws.Range("J92").Value = 387
$ws.Range("K92").Value = 260.05
$ws.Range("L92").Value = 387
$ws.Range("M92").Value = 987.95
$ws.Range("N92").Value = -2883
$ws.Range("H93").Value = 98436.5
$ws.Range("J93").Value = 98436.5
$ws.Range("L93").Value = 98436.5
$ws.Range("N93").Value = -103428.5
$ws.Range("H112").Value = 1206
$ws.Range("I112").Value = 1166.6666
$ws.Range("J112").Value = 1210.3704
$ws.Range("K112").Value = 3499.9998
$ws.Range("L112").Value = 3631.1112
$ws.Range("M112").Value = -2391.9998
$ws.Range("N112").Value = -5847.1112
$ws.Range("H118").Value = 945.6111
$ws.Range("I118").Value = 428.66666
$ws.Range("J118").Value = 1204.0834
$ws.Range("K118").Value = 1285.99998
$ws.Range("L118").Value = 3612.2502
$ws.Range("M118").Value = 371.0000199999999
$ws.Range("N118").Value = -6926.2502
$ws.Range("H125").Value = 841567.8
$ws.Range("I125").Value = 1367.4286
$ws.Range("J125").Value = 1681768.1
$ws.Range("K125").Value = 12306.8574
$ws.Range("L125").Value = 15135912.9
$ws.Range("M125").Value = -9846.857399999999
$ws.Range("N125").Value = -15140832.9
$ws.Range("H129").Value = 1097.9584
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 1097.9584
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 3293.8752
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -13293.8752
$ws.Range("H132").Value = 23932.666
$ws.Range("I132").Value = 26198.977
$ws.Range("J132").Value = 703
$ws.Range("K132").Value = 78596.931
$ws.Range("L132").Value = 2109
$ws.Range("M132").Value = -76066.931
$ws.Range("N132").Value = -7169
$ws.Range("H137").Value = 1341.7
$ws.Range("I137").Value = 1284.1904
$ws.Range("K137").Value = 3852.5712
$ws.Range("M137").Value = -1302.5712

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4793.1797
$ws.Range("I32").Value = 3401.378
$ws.Range("J32").Value = 21097.143
$ws.Range("K32").Value = 3401.378
$ws.Range("L32").Value = 21097.143
$ws.Range("M32").Value = -3114.378
$ws.Range("N32").Value = -21671.143
$ws.Range("H45").Value = 1340.909
$ws.Range("I45").Value = 1195.6666
$ws.Range("K45").Value = 1195.6666
$ws.Range("M45").Value = -818.6666
$ws.Range("H61").Value = 1357.2727
$ws.Range("I61").Value = 760.2222
$ws.Range("J61").Value = 2305.5293
$ws.Range("K61").Value = 760.2222
$ws.Range("L61").Value = 2305.5293
$ws.Range("M61").Value = -548.2222
$ws.Range("N61").Value = -2729.5293
$ws.Range("H74").Value = 3890.05
$ws.Range("I74").Value = 4155.222
$ws.Range("K74").Value = 4155.222
$ws.Range("M74").Value = -3281.222
$ws.Range("H77").Value = 3890.05
$ws.Range("I77").Value = 4155.222
$ws.Range("K77").Value = 20776.11
$ws.Range("M77").Value = -16408.11
$ws.Range("H122").Value = 1423.5454
$ws.Range("I122").Value = 1067.3462
$ws.Range("J122").Value = 2746.5715
$ws.Range("K122").Value = 3202.0386
$ws.Range("L122").Value = 8239.7145
$ws.Range("M122").Value = -752.0385999999999
$ws.Range("N122").Value = -13139.7145
$ws.Range("H132").Value = 1850.82
$ws.Range("I132").Value = 943.3461
$ws.Range("J132").Value = 2833.9167
$ws.Range("K132").Value = 2830.0383
$ws.Range("L132").Value = 8501.750100000001
$ws.Range("M132").Value = -300.0383000000002
$ws.Range("N132").Value = -13561.7501
$ws.Range("H136").Value = 1357.2727
$ws.Range("I136").Value = 760.2222
$ws.Range("J136").Value = 2305.5293
$ws.Range("K136").Value = 2280.6666
$ws.Range("L136").Value = 6916.5879
$ws.Range("M136").Value = 269.3334
$ws.Range("N136").Value = -12016.5879

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 650.7368
$ws.Range("I94").Value = 685.5294
$ws.Range("J94").Value = 355
$ws.Range("K94").Value = 685.5294
$ws.Range("L94").Value = 355
$ws.Range("M94").Value = -234.5294
$ws.Range("N94").Value = -1257

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2525.2646
$ws.Range("I31").Value = 1608.3823
$ws.Range("J31").Value = 3442.147
$ws.Range("K31").Value = 1608.3823
$ws.Range("L31").Value = 3442.147
$ws.Range("M31").Value = -1313.3823
$ws.Range("N31").Value = -4032.147
$ws.Range("H34").Value = 2525.2646
$ws.Range("I34").Value = 1608.3823
$ws.Range("J34").Value = 3442.147
$ws.Range("K34").Value = 1608.3823
$ws.Range("L34").Value = 3442.147
$ws.Range("M34").Value = -1406.3823
$ws.Range("N34").Value = -3846.147
$ws.Range("H132").Value = 1616.7805
$ws.Range("I132").Value = 810.4583
$ws.Range("J132").Value = 2755.1177
$ws.Range("K132").Value = 2431.3749
$ws.Range("L132").Value = 8265.3531
$ws.Range("M132").Value = 98.6251000000002
$ws.Range("N132").Value = -13325.3531

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 878637.7
$ws.Range("I5").Value = 1135.5
$ws.Range("J5").Value = 1463639.1
$ws.Range("K5").Value = 3406.5
$ws.Range("L5").Value = 4390917.300000001
$ws.Range("M5").Value = -3294.5
$ws.Range("N5").Value = -4391141.300000001
$ws.Range("H113").Value = 1136.0555
$ws.Range("I113").Value = 1169.9333
$ws.Range("J113").Value = 966.6667
$ws.Range("K113").Value = 3509.7999
$ws.Range("L113").Value = 2900.0001
$ws.Range("M113").Value = -1339.7999
$ws.Range("N113").Value = -7240.0001
$ws.Range("H122").Value = 1165.96
$ws.Range("J122").Value = 2081.7273
$ws.Range("L122").Value = 18735.5457
$ws.Range("N122").Value = -23635.5457
$ws.Range("H132").Value = 1632.5
$ws.Range("I132").Value = 1828.8572
$ws.Range("J132").Value = 1560.1578
$ws.Range("K132").Value = 16459.7148
$ws.Range("L132").Value = 14041.4202
$ws.Range("M132").Value = -13929.7148
$ws.Range("N132").Value = -19101.4202
$ws.Range("H135").Value = 878637.7
$ws.Range("I135").Value = 1135.5
$ws.Range("J135").Value = 1463639.1
$ws.Range("K135").Value = 10219.5
$ws.Range("L135").Value = 13172751.9
$ws.Range("M135").Value = -7684.5
$ws.Range("N135").Value = -13177821.9
$ws.Range("H139").Value = 1645.862
$ws.Range("I139").Value = 1133.0769
$ws.Range("J139").Value = 2062.5
$ws.Range("K139").Value = 3399.2307
$ws.Range("L139").Value = 6187.5
$ws.Range("M139").Value = 1740.7693
$ws.Range("N139").Value = -16467.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2146.75
$ws.Range("I122").Value = 1860.1875
$ws.Range("J122").Value = 2528.8333
$ws.Range("K122").Value = 5580.5625
$ws.Range("L122").Value = 7586.499899999999
$ws.Range("M122").Value = -3130.5625
$ws.Range("N122").Value = -12486.4999

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 79931.08
$ws.Range("I7").Value = 126763
$ws.Range("K7").Value = 126763
$ws.Range("M7").Value = -126651
$ws.Range("H40").Value = 23732.84
$ws.Range("I40").Value = 29137.947
$ws.Range("K40").Value = 29137.947
$ws.Range("M40").Value = -29001.947
$ws.Range("H55").Value = 564.8
$ws.Range("I55").Value = 309.2
$ws.Range("J55").Value = 820.4
$ws.Range("K55").Value = 309.2
$ws.Range("L55").Value = 820.4
$ws.Range("M55").Value = -136.2
$ws.Range("N55").Value = -1166.4
$ws.Range("H82").Value = 3160.7
$ws.Range("I82").Value = 1178.5
$ws.Range("J82").Value = 4482.1665
$ws.Range("K82").Value = 1178.5
$ws.Range("L82").Value = 4482.1665
$ws.Range("M82").Value = -817.5
$ws.Range("N82").Value = -5204.1665
$ws.Range("H85").Value = 3160.7
$ws.Range("I85").Value = 1178.5
$ws.Range("J85").Value = 4482.1665
$ws.Range("K85").Value = 1178.5
$ws.Range("L85").Value = 4482.1665
$ws.Range("M85").Value = 69.5
$ws.Range("N85").Value = -6978.1665
$ws.Range("H126").Value = 79931.08
$ws.Range("I126").Value = 126763
$ws.Range("K126").Value = 380289
$ws.Range("M126").Value = -377819
$ws.Range("H132").Value = 4866.2856
$ws.Range("I132").Value = 4433.797
$ws.Range("J132").Value = 6250.25
$ws.Range("K132").Value = 13301.391
$ws.Range("L132").Value = 18750.75
$ws.Range("M132").Value = -10771.391
$ws.Range("N132").Value = -23810.75
$ws.Range("H136").Value = 11496176
$ws.Range("I136").Value = 2105.9614
$ws.Range("J136").Value = 111111460
$ws.Range("K136").Value = 6317.8842
$ws.Range("L136").Value = 333334380
$ws.Range("M136").Value = -3767.8842
$ws.Range("N136").Value = -333339480

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 891.88
$ws.Range("I113").Value = 974.6111
$ws.Range("J113").Value = 679.1429000000001
$ws.Range("K113").Value = 2923.8333
$ws.Range("L113").Value = 2037.4287
$ws.Range("M113").Value = -753.8332999999998
$ws.Range("N113").Value = -6377.4287
$ws.Range("H126").Value = 26504.385
$ws.Range("I126").Value = 28560.305
$ws.Range("J126").Value = 1833.3334
$ws.Range("K126").Value = 85680.91500000001
$ws.Range("L126").Value = 5500.0002
$ws.Range("M126").Value = -83210.91500000001
$ws.Range("N126").Value = -10440.0002
$ws.Range("H132").Value = 1608.8334
$ws.Range("I132").Value = 1472.8334
$ws.Range("K132").Value = 4418.5002
$ws.Range("M132").Value = -1888.5002
$ws.Range("H136").Value = 5651312
$ws.Range("I136").Value = 8333911
$ws.Range("J136").Value = 3735.2632
$ws.Range("K136").Value = 25001733
$ws.Range("L136").Value = 11205.7896
$ws.Range("M136").Value = -24999183
$ws.Range("N136").Value = -16305.7896
